$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update F column values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12273
$ws1.Range("F3").Value = 57
$ws1.Range("F4").Value = 247
$ws1.Range("F7").Value = 12209
$ws1.Range("F8").Value = 515
$ws1.Range("F10").Value = 119
$ws1.Range("F11").Value = 618
$ws1.Range("F12").Value = 2811
$ws1.Range("F13").Value = 5985
$ws1.Range("F14").Value = 140
$ws1.Range("F15").Value = 3571

# Sheet "演出" (Performance) - G2 from 499 (number) to "不可售" (text)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value = "不可售"

# Sheet "全部类型" (All types) - G2 and F column updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F3").Value = 12273
$ws4.Range("F4").Value = 57
$ws4.Range("F5").Value = 247
$ws4.Range("F9").Value = 12209
$ws4.Range("F10").Value = 515
$ws4.Range("F12").Value = 119
$ws4.Range("F13").Value = 618
$ws4.Range("F14").Value = 2811
$ws4.Range("F16").Value = 5985
$ws4.Range("F17").Value = 140
$ws4.Range("F18").Value = 3571
